# Cập nhật report lương tổng hợp cho NV-9 Lê Văn Linh (tháng 8-2024)
# Sheet "Lương" chứa các giá trị tổng hợp tính theo công thực tế trong tháng.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# Tổng công tại CẦN THƠ: 8 -> 11 ngày công
$ws.Range("B2").Value = 11

# Phụ cấp tại CẦN THƠ (đơn giá 35.000đ / công)
$ws.Range("B3").Value = 385000

# Lương cơ bản tại CẦN THƠ
$ws.Range("B4").Value = 5892857.142857143

# Lương cơ bản tại LONG XUYÊN
$ws.Range("B15").Value = 3928571.428571429

# Lương cơ bản tại SÓC TRĂNG
$ws.Range("B26").Value = 5892857.142857143

# Tổng lương tại CẦN THƠ
$ws.Range("B35").Value = 7477857.142857143

# Tổng lương tại LONG XUYÊN
$ws.Range("B36").Value = 3928571.428571429

# Tổng lương tại SÓC TRĂNG
$ws.Range("B37").Value = 5892857.142857143

# Tổng lương tại HỆ THỐNG
$ws.Range("B38").Value = 17299285.71428571
